$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" contain identical data tables that both need
# their "想去人数" (column F) figures refreshed for rows 4-6.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 164
    $ws.Range("F5").Value = 83
    $ws.Range("F6").Value = 699
}
